# Update "想去人数" (want-to-go count) figures for two events that are
# duplicated across the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 195
    $ws.Range("F4").Value = 144
}
